# Add two new columns, I ("I0") and J ("IF"), to the single-sheet workbook,
# extending the used range from A1:H43 to A1:J43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells (I1, J1) the same formatting (bold, centered,
# thin-bordered) as the existing header row by copying the style from H1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-43 in the new I0 / IF columns
$iVals = @(8,8,8,8,8,8,7,9,6,6,8,10,8,5,7,6,7,1,7,9,5,4,8,7,6,8,6,9,4,7,4,4,9,9,7,6,3,6,5,4,6,6)
$jVals = @(8,8,8,8,8,8,7,9,7,7,8,10,8,5,7,6,7,2,7,9,6,5,8,8,7,8,6,9,6,7,4,4,9,9,7,6,3,6,6,4,6,6)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
